$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value2 = 3
$ws.Range("F2").Value2 = 1
$ws.Range("G2").Value2 = 37.63904266666666
$ws.Range("H2").Value2 = 112.917128
$ws.Range("I2").Value2 = 0.4850220755088102
$ws.Range("J2").Value2 = 0.4850220755088102
$ws.Range("K2").Value2 = 3
$ws.Range("L2").Value2 = 1
$ws.Range("M2").Value2 = 15.35884066666667
$ws.Range("N2").Value2 = 46.076522
$ws.Range("O2").Value2 = 0.1012042817263867
$ws.Range("P2").Value2 = 0.1012042817263867
$ws.Range("Q2").Value2 = 578.0920591632016
$ws.Range("R2").Value2 = 5202.828532468815
$ws.Range("S2").Value2 = 0.0490863107733104
$ws.Range("T2").Value2 = 0.04908631077331042

# Row 3
$ws.Range("E3").Value2 = 3
$ws.Range("F3").Value2 = 1
$ws.Range("G3").Value2 = 37.63904266666666
$ws.Range("H3").Value2 = 112.917128
$ws.Range("I3").Value2 = 0.4850220755088102
$ws.Range("J3").Value2 = 0.4850220755088102
$ws.Range("K3").Value2 = 3
$ws.Range("L3").Value2 = 1
$ws.Range("M3").Value2 = 50.59256466666667
$ws.Range("N3").Value2 = 151.777694
$ws.Range("O3").Value2 = 0.3333704853712116
$ws.Range("P3").Value2 = 0.3333704853712116
$ws.Range("Q3").Value2 = 1904.255700104759
$ws.Range("R3").Value2 = 17138.30130094283
$ws.Range("S3").Value2 = 0.1616920447281245
$ws.Range("T3").Value2 = 0.1616920447281245

# Row 4
$ws.Range("E4").Value2 = 3
$ws.Range("F4").Value2 = 1
$ws.Range("G4").Value2 = 37.63904266666666
$ws.Range("H4").Value2 = 112.917128
$ws.Range("I4").Value2 = 0.4850220755088102
$ws.Range("J4").Value2 = 0.4850220755088102
$ws.Range("K4").Value2 = 3
$ws.Range("L4").Value2 = 1
$ws.Range("M4").Value2 = 60.37715666666667
$ws.Range("N4").Value2 = 181.13147
$ws.Range("O4").Value2 = 0.397844271305776
$ws.Range("P4").Value2 = 0.397844271305776
$ws.Range("Q4").Value2 = 2272.538375868684
$ws.Range("R4").Value2 = 20452.84538281816
$ws.Range("S4").Value2 = 0.1929632541980176
$ws.Range("T4").Value2 = 0.1929632541980177

# Row 5
$ws.Range("E5").Value2 = 3
$ws.Range("F5").Value2 = 1
$ws.Range("G5").Value2 = 37.63904266666666
$ws.Range("H5").Value2 = 112.917128
$ws.Range("I5").Value2 = 0.4850220755088102
$ws.Range("J5").Value2 = 0.4850220755088102
$ws.Range("K5").Value2 = 3
$ws.Range("L5").Value2 = 1
$ws.Range("M5").Value2 = 25.43221733333333
$ws.Range("N5").Value2 = 76.29665199999999
$ws.Range("O5").Value2 = 0.1675809615966257
$ws.Range("P5").Value2 = 0.1675809615966258
$ws.Range("Q5").Value2 = 957.2443133172726
$ws.Range("R5").Value2 = 8615.198819855455
$ws.Range("S5").Value2 = 0.08128046580935763
$ws.Range("T5").Value2 = 0.08128046580935765

# Row 6
$ws.Range("E6").Value2 = 3
$ws.Range("F6").Value2 = 1
$ws.Range("G6").Value2 = 17.57434666666667
$ws.Range("H6").Value2 = 52.72304
$ws.Range("I6").Value2 = 0.2264655392929762
$ws.Range("J6").Value2 = 0.2264655392929762
$ws.Range("K6").Value2 = 3
$ws.Range("L6").Value2 = 1
$ws.Range("M6").Value2 = 15.35884066666667
$ws.Range("N6").Value2 = 46.076522
$ws.Range("O6").Value2 = 0.1012042817263867
$ws.Range("P6").Value2 = 0.1012042817263867
$ws.Range("Q6").Value2 = 269.9215902740978
$ws.Range("R6").Value2 = 2429.29431246688
$ws.Range("S6").Value2 = 0.02291928223992445
$ws.Range("T6").Value2 = 0.02291928223992445

# Row 7
$ws.Range("E7").Value2 = 3
$ws.Range("F7").Value2 = 1
$ws.Range("G7").Value2 = 17.57434666666667
$ws.Range("H7").Value2 = 52.72304
$ws.Range("I7").Value2 = 0.2264655392929762
$ws.Range("J7").Value2 = 0.2264655392929762
$ws.Range("K7").Value2 = 3
$ws.Range("L7").Value2 = 1
$ws.Range("M7").Value2 = 50.59256466666667
$ws.Range("N7").Value2 = 151.777694
$ws.Range("O7").Value2 = 0.3333704853712116
$ws.Range("P7").Value2 = 0.3333704853712116
$ws.Range("Q7").Value2 = 889.1312702077512
$ws.Range("R7").Value2 = 8002.18143186976
$ws.Range("S7").Value2 = 0.07549692675395266
$ws.Range("T7").Value2 = 0.07549692675395266

# Row 8
$ws.Range("E8").Value2 = 3
$ws.Range("F8").Value2 = 1
$ws.Range("G8").Value2 = 17.57434666666667
$ws.Range("H8").Value2 = 52.72304
$ws.Range("I8").Value2 = 0.2264655392929762
$ws.Range("J8").Value2 = 0.2264655392929762
$ws.Range("K8").Value2 = 3
$ws.Range("L8").Value2 = 1
$ws.Range("M8").Value2 = 60.37715666666667
$ws.Range("N8").Value2 = 181.13147
$ws.Range("O8").Value2 = 0.397844271305776
$ws.Range("P8").Value2 = 0.397844271305776
$ws.Range("Q8").Value2 = 1061.089082007644
$ws.Range("R8").Value2 = 9549.8017380688
$ws.Range("S8").Value2 = 0.09009801745588369
$ws.Range("T8").Value2 = 0.09009801745588369

# Row 9
$ws.Range("E9").Value2 = 3
$ws.Range("F9").Value2 = 1
$ws.Range("G9").Value2 = 17.57434666666667
$ws.Range("H9").Value2 = 52.72304
$ws.Range("I9").Value2 = 0.2264655392929762
$ws.Range("J9").Value2 = 0.2264655392929762
$ws.Range("K9").Value2 = 3
$ws.Range("L9").Value2 = 1
$ws.Range("M9").Value2 = 25.43221733333333
$ws.Range("N9").Value2 = 76.29665199999999
$ws.Range("O9").Value2 = 0.1675809615966257
$ws.Range("P9").Value2 = 0.1675809615966258
$ws.Range("Q9").Value2 = 446.9546039180088
$ws.Range("R9").Value2 = 4022.591435262079
$ws.Range("S9").Value2 = 0.03795131284321538
$ws.Range("T9").Value2 = 0.03795131284321539

# Row 10
$ws.Range("E10").Value2 = 3
$ws.Range("F10").Value2 = 1
$ws.Range("G10").Value2 = 19.168158
$ws.Range("H10").Value2 = 57.504474
$ws.Range("I10").Value2 = 0.2470036195972184
$ws.Range("J10").Value2 = 0.2470036195972184
$ws.Range("K10").Value2 = 3
$ws.Range("L10").Value2 = 1
$ws.Range("M10").Value2 = 15.35884066666667
$ws.Range("N10").Value2 = 46.076522
$ws.Range("O10").Value2 = 0.1012042817263867
$ws.Range("P10").Value2 = 0.1012042817263867
$ws.Range("Q10").Value2 = 294.400684595492
$ws.Range("R10").Value2 = 2649.606161359428
$ws.Range("S10").Value2 = 0.02499782390515413
$ws.Range("T10").Value2 = 0.02499782390515413

# Row 11
$ws.Range("E11").Value2 = 3
$ws.Range("F11").Value2 = 1
$ws.Range("G11").Value2 = 19.168158
$ws.Range("H11").Value2 = 57.504474
$ws.Range("I11").Value2 = 0.2470036195972184
$ws.Range("J11").Value2 = 0.2470036195972184
$ws.Range("K11").Value2 = 3
$ws.Range("L11").Value2 = 1
$ws.Range("M11").Value2 = 50.59256466666667
$ws.Range("N11").Value2 = 151.777694
$ws.Range("O11").Value2 = 0.3333704853712116
$ws.Range("P11").Value2 = 0.3333704853712116
$ws.Range("Q11").Value2 = 969.7662731558842
$ws.Range("R11").Value2 = 8727.896458402956
$ws.Range("S11").Value2 = 0.0823437165535708
$ws.Range("T11").Value2 = 0.0823437165535708

# Row 12
$ws.Range("E12").Value2 = 3
$ws.Range("F12").Value2 = 1
$ws.Range("G12").Value2 = 19.168158
$ws.Range("H12").Value2 = 57.504474
$ws.Range("I12").Value2 = 0.2470036195972184
$ws.Range("J12").Value2 = 0.2470036195972184
$ws.Range("K12").Value2 = 3
$ws.Range("L12").Value2 = 1
$ws.Range("M12").Value2 = 60.37715666666667
$ws.Range("N12").Value2 = 181.13147
$ws.Range("O12").Value2 = 0.397844271305776
$ws.Range("P12").Value2 = 0.397844271305776
$ws.Range("Q12").Value2 = 1157.31887857742
$ws.Range("R12").Value2 = 10415.86990719678
$ws.Range("S12").Value2 = 0.09826897504854444
$ws.Range("T12").Value2 = 0.09826897504854444

# Row 13
$ws.Range("E13").Value2 = 3
$ws.Range("F13").Value2 = 1
$ws.Range("G13").Value2 = 19.168158
$ws.Range("H13").Value2 = 57.504474
$ws.Range("I13").Value2 = 0.2470036195972184
$ws.Range("J13").Value2 = 0.2470036195972184
$ws.Range("K13").Value2 = 3
$ws.Range("L13").Value2 = 1
$ws.Range("M13").Value2 = 25.43221733333333
$ws.Range("N13").Value2 = 76.29665199999999
$ws.Range("O13").Value2 = 0.1675809615966257
$ws.Range("P13").Value2 = 0.1675809615966258
$ws.Range("Q13").Value2 = 487.488760135672
$ws.Range("R13").Value2 = 4387.398841221047
$ws.Range("S13").Value2 = 0.04139310408994901
$ws.Range("T13").Value2 = 0.04139310408994902

# Row 14
$ws.Range("E14").Value2 = 3
$ws.Range("F14").Value2 = 1
$ws.Range("G14").Value2 = 3.221194
$ws.Range("H14").Value2 = 9.663582
$ws.Range("I14").Value2 = 0.04150876560099527
$ws.Range("J14").Value2 = 0.04150876560099527
$ws.Range("K14").Value2 = 3
$ws.Range("L14").Value2 = 1
$ws.Range("M14").Value2 = 15.35884066666667
$ws.Range("N14").Value2 = 46.076522
$ws.Range("O14").Value2 = 0.1012042817263867
$ws.Range("P14").Value2 = 0.1012042817263867
$ws.Range("Q14").Value2 = 49.47380540242266
$ws.Range("R14").Value2 = 445.264248621804
$ws.Range("S14").Value2 = 0.004200864807997673
$ws.Range("T14").Value2 = 0.004200864807997673

# Row 15
$ws.Range("E15").Value2 = 3
$ws.Range("F15").Value2 = 1
$ws.Range("G15").Value2 = 3.221194
$ws.Range("H15").Value2 = 9.663582
$ws.Range("I15").Value2 = 0.04150876560099527
$ws.Range("J15").Value2 = 0.04150876560099527
$ws.Range("K15").Value2 = 3
$ws.Range("L15").Value2 = 1
$ws.Range("M15").Value2 = 50.59256466666667
$ws.Range("N15").Value2 = 151.777694
$ws.Range("O15").Value2 = 0.3333704853712116
$ws.Range("P15").Value2 = 0.3333704853712116
$ws.Range("Q15").Value2 = 162.9684657488787
$ws.Range("R15").Value2 = 1466.716191739908
$ws.Range("S15").Value2 = 0.01383779733556364
$ws.Range("T15").Value2 = 0.01383779733556364

# Row 16
$ws.Range("E16").Value2 = 3
$ws.Range("F16").Value2 = 1
$ws.Range("G16").Value2 = 3.221194
$ws.Range("H16").Value2 = 9.663582
$ws.Range("I16").Value2 = 0.04150876560099527
$ws.Range("J16").Value2 = 0.04150876560099527
$ws.Range("K16").Value2 = 3
$ws.Range("L16").Value2 = 1
$ws.Range("M16").Value2 = 60.37715666666667
$ws.Range("N16").Value2 = 181.13147
$ws.Range("O16").Value2 = 0.397844271305776
$ws.Range("P16").Value2 = 0.397844271305776
$ws.Range("Q16").Value2 = 194.4865347917267
$ws.Range("R16").Value2 = 1750.37881312554
$ws.Range("S16").Value2 = 0.01651402460333022
$ws.Range("T16").Value2 = 0.01651402460333022

# Row 17
$ws.Range("E17").Value2 = 3
$ws.Range("F17").Value2 = 1
$ws.Range("G17").Value2 = 3.221194
$ws.Range("H17").Value2 = 9.663582
$ws.Range("I17").Value2 = 0.04150876560099527
$ws.Range("J17").Value2 = 0.04150876560099527
$ws.Range("K17").Value2 = 3
$ws.Range("L17").Value2 = 1
$ws.Range("M17").Value2 = 25.43221733333333
$ws.Range("N17").Value2 = 76.29665199999999
$ws.Range("O17").Value2 = 0.1675809615966257
$ws.Range("P17").Value2 = 0.1675809615966258
$ws.Range("Q17").Value2 = 81.92210588082932
$ws.Range("R17").Value2 = 737.2989529274639
$ws.Range("S17").Value2 = 0.006956078854103728
$ws.Range("T17").Value2 = 0.00695607885410373
